$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment"); this shifts C:Q left to B:P.
$ws.Columns("B").Delete()

# Append ".deja.deja.deja" to every header in row 1 except "Country" (A1).
for ($col = 2; $col -le 16; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Value2 + ".deja.deja.deja"
}
